$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 300, shifting the existing
# rows 300-366 down to 302-368 (dimension grows from R366 to R368).
$ws.Range("A300:R301").Insert()

# New record for (new) row 300
$ws.Range("A300").Value = 6
$ws.Range("B300").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C300").Value = "Metropolitana"
$ws.Range("D300").Value = 45209
$ws.Range("E300").Value = 13
$ws.Range("F300").Value = 100112001
$ws.Range("G300").Value = "Berenjena"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 290
$ws.Range("K300").Value = 10000
$ws.Range("L300").Value = 12000
$ws.Range("M300").Value = 11172
$ws.Range("N300").Value = "$/caja 40 unidades"
$ws.Range("O300").Value = "Provincia de Huasco"
$ws.Range("P300").Value = 279
$ws.Range("Q300").Value = 40
$ws.Range("R300").Value = "Hortaliza"

# New record for (new) row 301
$ws.Range("A301").Value = 6
$ws.Range("B301").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C301").Value = "Metropolitana"
$ws.Range("D301").Value = 45209
$ws.Range("E301").Value = 13
$ws.Range("F301").Value = 100112001
$ws.Range("G301").Value = "Berenjena"
$ws.Range("H301").Value = "Sin especificar"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 260
$ws.Range("K301").Value = 6000
$ws.Range("L301").Value = 6500
$ws.Range("M301").Value = 6269
$ws.Range("N301").Value = "$/caja 50 unidades"
$ws.Range("O301").Value = "Región de Arica y Parinacota"
$ws.Range("P301").Value = 125
$ws.Range("Q301").Value = 50
$ws.Range("R301").Value = "Hortaliza"
